$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (NumberFormat "@") for Price cells whose new values would
# otherwise be auto-coerced to numbers by Excel, so they stay literal text like the source data.
$textForceCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D14", "D15", "D17", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D51")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated cell values (prices, volumes, and a few re-ranked coin rows).
$ws.Range("D2").Value = "36.551.89"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.104.52"
$ws.Range("E3").Value = "  +9.54%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "252.09"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "0.660"
$ws.Range("E6").Value = "  -6.22%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "47.43"
$ws.Range("E8").Value = "  +7.13%  "
$ws.Range("D9").Value = "59.27"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "0.372"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").Value = "0.0744"
$ws.Range("E11").Value = "  -3.02%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "2.414.47"
$ws.Range("D14").Value = "14.25"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "0.826"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "2.109.18"
$ws.Range("E16").Value = "  +9.88%  "
$ws.Range("D17").Value = "5.08"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "36.561.40"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "72.86"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").Value = "0.0₃0829"
$ws.Range("E20").Value = "  -4.08%  "
$ws.Range("D21").Value = "13.21"
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("D22").Value = "239.67"
$ws.Range("E22").Value = "  -4.64%  "
$ws.Range("D23").Value = "5.15"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  -8.77%  "
$ws.Range("D26").Value = "171.40"
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("D27").Value = "21.33"
$ws.Range("E27").Value = "  +13.69%  "
$ws.Range("D28").Value = "9.14"
$ws.Range("E28").Value = "  +2.94%  "
$ws.Range("D29").Value = "1.97"
$ws.Range("E29").Value = "  -10.84%  "
$ws.Range("D30").Value = "28.04"
$ws.Range("E30").Value = "  +59.47%  "
$ws.Range("E31").Value = "  -5.20%  "
$ws.Range("D32").Value = "4.43"
$ws.Range("E32").Value = "  -4.97%  "
$ws.Range("D33").Value = "0.0606"
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("D34").Value = "0.0886"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.942"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.32"
$ws.Range("E38").Value = "  +11.99%  "
$ws.Range("D39").Value = "4.06"
$ws.Range("E39").Value = "  -6.92%  "
$ws.Range("D40").Value = "1.33"
$ws.Range("E40").Value = "  -14.46%  "
$ws.Range("E41").Value = "  +5.56%  "
$ws.Range("D42").Value = "0.0222"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").Value = "97.66"
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "2.75"
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "16.22"
$ws.Range("E45").Value = "  -7.64%  "
$ws.Range("D46").Value = "1.333.26"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "0.0840"
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("D48").Value = "7.04"
$ws.Range("E48").Value = "  +9.45%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "2.84"
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.280.34"
$ws.Range("E50").Value = "  +8.13%  "
$ws.Range("D51").Value = "2.24"
$ws.Range("E51").Value = "  -6.85%  "
